$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Launch Parameters")

# Update the Turn Altitude[km] column (C) for every vessel from 0.25 to 0.5
$ws.Range("C2:C9").Value = 0.5

# Adjust Saturn V (row 5) specific launch parameters
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 40

# Update the selection shown in the sheet view
$ws.Range("C2:C9").Select()
$excel.ActiveCell = $ws.Range("C2")
